$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats
$xlVAlignCenter = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

# --- QOL: give the two new (currently blank) cells in row 2 the same
#     formatting as their neighbours (D2) ---
$ws.Range("D2").Copy()
$ws.Range("E2").PasteSpecial($xlPasteFormats)
$ws.Range("F2").PasteSpecial($xlPasteFormats)

# --- New row: add the JapaneseSamurai dummy alongside the existing Dummy/RubyRose rows ---
# Column A (name) - copy formatting from the Dummy name cell, then set the text
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial($xlPasteFormats)
$ws.Range("A4").Value = "JapaneseSamurai"

# Column B (skill) - stays blank, just inherit formatting
$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial($xlPasteFormats)

# Columns C-F (HP, aspd, atk, mana) - inherit formatting from row 3 and set values
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial($xlPasteFormats)
$ws.Range("C4").Value = 50

$ws.Range("D3").Copy()
$ws.Range("D4").PasteSpecial($xlPasteFormats)
$ws.Range("D4").Value = 1

$ws.Range("E3").Copy()
$ws.Range("E4").PasteSpecial($xlPasteFormats)
$ws.Range("E4").Value = 3.8

$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial($xlPasteFormats)
$ws.Range("F4").Value = 100

# --- Behaviour tweak: vertically center the long-range skill description text
#     (it was top-aligned before) so it reads better for the new QOL pass ---
$ws.Range("B3").VerticalAlignment = $xlVAlignCenter

# --- QOL: widen column A a bit now that it holds longer character names ---
$ws.Columns.Item(1).ColumnWidth = 14

# --- Leave the selection where the editor was last working ---
[void]$ws.Range("H9").Select()
